$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 176, pushing the existing
# rows 176-180 down to 178-182 (matches the diff: new weekly entries are
# added above the older rows, which keep their data but shift down by 2).
$ws.Rows("176:177").Insert()

# Column D uses a custom date/time number format (style index 2 in the
# original file) - copy it onto the two freshly-inserted date cells.
$dateFormat = $ws.Cells.Item(178, 4).NumberFormat
$ws.Cells.Item(176, 4).NumberFormat = $dateFormat
$ws.Cells.Item(177, 4).NumberFormat = $dateFormat

# Row 176 (new): Albahaca, Primera, Region de Arica y Parinacota, $/paquete
$ws.Cells.Item(176, 1).Value = 9
$ws.Cells.Item(176, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(176, 3).Value = "Metropolitana"
$ws.Cells.Item(176, 4).Value = 44448
$ws.Cells.Item(176, 5).Value = 13
$ws.Cells.Item(176, 6).Value = 100112052
$ws.Cells.Item(176, 7).Value = "Albahaca"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 169
$ws.Cells.Item(176, 11).Value = 5000
$ws.Cells.Item(176, 12).Value = 5500
$ws.Cells.Item(176, 13).Value = 5251
$ws.Cells.Item(176, 14).Value = "$/paquete"
$ws.Cells.Item(176, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(176, 16).Value = 5251
$ws.Cells.Item(176, 17).Value = 1
$ws.Cells.Item(176, 18).Value = "Hortaliza"

# Row 177 (new): Albahaca, Segunda, Region de Arica y Parinacota, $/paquete
$ws.Cells.Item(177, 1).Value = 9
$ws.Cells.Item(177, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(177, 3).Value = "Metropolitana"
$ws.Cells.Item(177, 4).Value = 44448
$ws.Cells.Item(177, 5).Value = 13
$ws.Cells.Item(177, 6).Value = 100112052
$ws.Cells.Item(177, 7).Value = "Albahaca"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Segunda"
$ws.Cells.Item(177, 10).Value = 97
$ws.Cells.Item(177, 11).Value = 4500
$ws.Cells.Item(177, 12).Value = 4500
$ws.Cells.Item(177, 13).Value = 4500
$ws.Cells.Item(177, 14).Value = "$/paquete"
$ws.Cells.Item(177, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(177, 16).Value = 4500
$ws.Cells.Item(177, 17).Value = 1
$ws.Cells.Item(177, 18).Value = "Hortaliza"
